# LKS.xlsx — add a new player ("A. Tutyškinas") as a new row inserted
# right before "M. Wszołek" (currently spreadsheet row 14).
#
# The sheet's column A holds a plain sequential index (0,1,2,...) that is
# NOT tied to the player on that line, so on insertion it is simply
# extended by one more value at the bottom (A32 = 30); the existing A14..A31
# index cells are left completely untouched. Every other column (B..N) for
# rows 14..31 is overwritten with the data that used to sit one row above
# it (i.e. the whole roster table shifts down by one row), row 14 gets the
# brand-new player's stats, and a new trailing row 32 is created that is a
# copy of the previous last row (player name "K. Moskal" with the rest of
# the stat columns left blank), matching the pre-existing row's shape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target contents for B14:N32 after the shift. First element of each inner
# array is the spreadsheet row number; the rest are the values for columns
# B..N in order. $null means "leave that cell blank" (no value written).
$rows = @(
  @(14, 'A. Tutyškinas', '19', 'O', '3', '1', '0', '1', '0', '1', '0', '0', '0', '0'),
  @(15, 'M. Wszołek', '20', 'O', '247', '4', '4', '0', '4', '6', '0', '0', '0', '0'),
  @(16, 'B. Biel', '28', 'P', '522', '10', '7', '3', '6', '3', '0', '1', '0', '0'),
  @(17, 'K. Ibe-Torti', '20', 'P', '341', '10', '2', '8', '2', '8', '0', '0', '0', '0'),
  @(18, 'D. Kort', '27', 'P', '536', '9', '7', '2', '7', '3', '0', '2', '0', '0'),
  @(19, 'M. Kowalczyk', '18', 'P', '602', '9', '7', '2', '4', '2', '3', '5', '0', '0'),
  @(20, 'J. Kuźma', '19', 'P', '257', '7', '2', '5', '2', '8', '0', '1', '0', '0'),
  @(21, 'Javi Moreno', '25', 'P', '81', '5', '0', '5', '0', '8', '0', '0', '0', '0'),
  @(22, 'D. Nowacki', '24', 'P', '0', '0', '0', '0', '0', '2', '0', '0', '0', '0'),
  @(23, 'V. Okhronchuk', '25', 'P', '167', '6', '2', '4', '2', '8', '0', '1', '0', '0'),
  @(24, 'J. Romanowicz', '21', 'P', '0', '0', '0', '0', '0', '0', '0', '0', '0', '0'),
  @(25, 'Pirulo', '30', 'P', '836', '10', '10', '0', '5', '0', '6', '2', '1', '0'),
  @(26, 'M. Trąbka', '25', 'P', '810', '10', '9', '1', '4', '1', '2', '1', '0', '0'),
  @(27, 'N. Balongo', '23', 'N', '625', '10', '8', '2', '6', '2', '0', '2', '0', '0'),
  @(28, 'G. Glapka', '19', 'N', '2', '1', '0', '1', '0', '1', '0', '0', '0', '0'),
  @(29, 'P. Janczukowicz', '22', 'N', '272', '8', '2', '6', '2', '7', '1', '1', '0', '0'),
  @(30, 'S. Jurić', '23', 'N', '0', '0', '0', '0', '0', '0', '0', '0', '0', '0'),
  @(31, 'M. Radaszkiewicz', '25', 'N', '20', '1', '0', '1', '0', '1', '0', '0', '0', '0'),
  @(32, 'K. Moskal', $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
)

foreach ($row in $rows) {
  $r = $row[0]

  # Column B (player name) is never numeric-looking, write it plainly.
  $ws.Cells.Item($r, 2).Value = $row[1]

  # Columns C..N hold numbers-as-text in this sheet (age, minutes played,
  # appearance counts, …). Force the Text number format before assigning
  # so a value like "19" is not silently reinterpreted as the number 19.
  for ($c = 3; $c -le 14; $c++) {
    $val = $row[$c - 1]
    if ($val -ne $null) {
      $cell = $ws.Cells.Item($r, $c)
      $cell.NumberFormat = "@"
      $cell.Value = $val
    }
  }
}

# New row 32 needs an index value in column A (continuing the existing
# 0-based sequence: row32 -> 30). Copy the format from A31 first so it
# reuses the existing bold/bordered/centered style instead of minting a
# new one, then set its numeric value.
$ws.Range("A31").Copy()
$ws.Range("A32").PasteSpecial(-4122)
$ws.Cells.Item(32, 1).Value = 30
